$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix death functionality bug: the "Shield Capacity" max-stat input (Q40)
# was wrong. Correcting it from 750 to 400 ripples through the calculated
# Table2 column (Q3:Q37) and every downstream ratio column (U:AB) via the
# worksheet's existing formulas.
$ws.Range("Q40").Value = 400

# Reflect where the author ended up scrolled to / selected while
# reviewing the fix.
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("Q41").Select()
